$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.237.11"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").Value = "1.819.69"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'313.03"
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("D7").Value = "'0.4449"
$ws.Range("E7").Value = "  -0.66%  "
$ws.Range("E8").Value = "  +1.84%  "
$ws.Range("D9").Value = "'0.07387"
$ws.Range("E9").Value = "  +1.51%  "
$ws.Range("D10").Value = "'0.8788"
$ws.Range("E10").Value = "  +2.86%  "
$ws.Range("E11").Value = "  +0.41%  "
$ws.Range("D12").Value = "1.819.93"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").Value = "'6.697"
$ws.Range("E13").Value = "  +0.87%  "
$ws.Range("D14").Value = "'5.416"
$ws.Range("E14").Value = "  +1.75%  "
$ws.Range("D15").Value = "'93.07"
$ws.Range("E15").Value = "  +0.76%  "
$ws.Range("D16").Value = "'0.07108"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").Value = "'1.002"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").Value = "'0.000008796"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").Value = "'15.00"
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("D21").Value = "27.254.57"
$ws.Range("E21").Value = "  +1.01%  "
$ws.Range("E22").Value = "  +3.69%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "'1.957"
$ws.Range("E24").Value = "  -1.38%  "
$ws.Range("D25").Value = "'151.04"
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("D26").Value = "'2.295"
$ws.Range("E26").Value = "  +3.04%  "
$ws.Range("D27").Value = "'18.56"
$ws.Range("E27").Value = "  +0.83%  "
$ws.Range("D28").Value = "'5.339"
$ws.Range("E28").Value = "  +1.90%  "
$ws.Range("D29").Value = "'117.23"
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("D30").Value = "'0.08863"
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("D31").Value = "'0.7865"
$ws.Range("E31").Value = "  +4.65%  "
$ws.Range("E32").Value = "  +0.91%  "
$ws.Range("D33").Value = "'4.556"
$ws.Range("E33").Value = "  +2.54%  "
$ws.Range("D34").Value = "'2.913"
$ws.Range("E34").Value = "  -1.80%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "'1.106"
$ws.Range("E36").Value = "  +1.09%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("E38").Value = "  +0.52%  "
$ws.Range("D39").Value = "'7.296"
$ws.Range("E39").Value = "  +1.92%  "
$ws.Range("D40").Value = "'0.5279"
$ws.Range("E40").Value = "  -0.70%  "
$ws.Range("D41").Value = "'2.864"
$ws.Range("E41").Value = "  -0.57%  "
$ws.Range("D42").Value = "'0.1700"
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("E43").Value = "  +15.91%  "
$ws.Range("E44").Value = "  +1.05%  "
$ws.Range("D45").Value = "'0.5026"
$ws.Range("E45").Value = "  -4.41%  "
$ws.Range("D46").Value = "'10.54"
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("D47").Value = "'104.86"
$ws.Range("E47").Value = "  -0.64%  "
$ws.Range("E48").Value = "  +1.21%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").Value = "'0.06382"
$ws.Range("D51").Value = "'65.98"
$ws.Range("E51").Value = "  +4.62%  "
